$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: variable names / labels
$ws.Cells.Item(1, 1).Value = "Intervalo renta"
$ws.Cells.Item(1, 2).Value = "Municipio código"
$ws.Cells.Item(1, 3).Value = "Año"
$ws.Cells.Item(1, 4).Value = "Municipio nombre"

# Row 2: sdmx / iaest identifiers
$ws.Cells.Item(2, 1).Value = "iaest-measure:intervalo-renta"
$ws.Cells.Item(2, 2).Value = "null"
$ws.Cells.Item(2, 3).Value = "sdmx-dimension:refPeriod"
$ws.Cells.Item(2, 4).Value = "sdmx-dimension:refArea"

# Row 3: dim / medida classifier
$ws.Cells.Item(3, 1).Value = "medida"
$ws.Cells.Item(3, 2).Value = "null"
$ws.Cells.Item(3, 3).Value = "dim"
$ws.Cells.Item(3, 4).Value = "dim"

# Row 4: datatype
$ws.Cells.Item(4, 1).Value = "xsd:int"
$ws.Cells.Item(4, 2).Value = "null"
$ws.Cells.Item(4, 3).Value = "xsd:date"
$ws.Cells.Item(4, 4).Value = "URI-Municipio"

# Row 5: mapping file moves from column A to column C, column A cleared
$ws.Cells.Item(5, 1).Clear()
$ws.Cells.Item(5, 3).Value = "mapping-ano.xlsx"
